# GEH Wards.xlsx - "Re-edit to remove Bed No. for MNH Pharmacy & Update Ward No"
#
# Ward 7W currently runs from Room 733 (row 127) through Room 764 (row 157),
# but rooms 737-764 were off by one (738 was listed twice/737 was skipped).
# Renumber rooms 738-764 (rows 131-157) down to 737-763, then add back the
# room that was missing (764) as a new row so 7W ends at 764 again, pushing
# every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Room No")

# Step 1: decrement the room numbers for rows 131-157 (738..764 -> 737..763)
for ($r = 131; $r -le 157; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 - 1
}

# Step 2: insert a new row at 158 (shifts old row 158.. down to 159..)
$ws.Rows.Item(158).Insert()

# Step 3: populate the newly inserted row with Room 764, Ward 7W
$ws.Cells.Item(158, 1).Value2 = 764
$ws.Cells.Item(158, 2).Value2 = "7W"

# Step 4: update the view/selection to where the author left off editing
$excel.ActiveWindow.ScrollRow = 96
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("E126").Select()
